{"js": "// Replace each old value with its new value, matching the commit diff.\n// Every old text value in this document is unique, so a simple\n// search-and-replace per pair is unambiguous and safe.\nconst replacements = [\n  [\"2025-04-14 Monday\", \"2025-04-15 Tuesday\"],\n  [\"870\u00d76=\", \"262\u00d75=\"],\n  [\"237\u00d76=\", \"771\u00d75=\"],\n  [\"356\u00d78=\", \"614\u00d76=\"],\n  [\"995\u00d72=\", \"138\u00d77=\"],\n  [\"534\u00d75=\", \"452\u00d77=\"],\n  [\"329\u00d77=\", \"806\u00d72=\"],\n  [\"196\u00d75=\", \"188\u00d78=\"],\n  [\"949\u00d79=\", \"394\u00d77=\"],\n  [\"921\u00d74=\", \"765\u00d73=\"],\n  [\"459\u00d72=\", \"832\u00d77=\"],\n  [\"242\u00d79=\", \"502\u00d75=\"],\n  [\"749\u00d72=\", \"553\u00d72=\"],\n  [\"692\u00d77=\", \"782\u00d75=\"],\n  [\"707\u00d78=\", \"718\u00d73=\"],\n  [\"315\u00d78=\", \"923\u00d78=\"],\n  [\"515\u00d76=\", \"892\u00d77=\"],\n  [\"398\u00d73=\", \"548\u00d77=\"],\n  [\"325\u00d77=\", \"975\u00d78=\"],\n  [\"465\u00d72=\", \"929\u00d78=\"],\n  [\"353\u00d78=\", \"641\u00d75=\"],\n  [\"647\u00d72=\", \"325\u00d76=\"],\n  [\"963\u00d72=\", \"336\u00d74=\"],\n  [\"422\u00d74=\", \"776\u00d77=\"],\n  [\"645\u00d79=\", \"586\u00d76=\"],\n  [\"362\u00d76=\", \"489\u00d75=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old value with its new value, matching the commit diff.\n# Every old text value in this document is unique, so Find/Replace per\n# pair (wdReplaceAll = 2, wdFindContinue = 1) is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-14 Monday\", \"2025-04-15 Tuesday\"),\n    @(\"870\u00d76=\", \"262\u00d75=\"),\n    @(\"237\u00d76=\", \"771\u00d75=\"),\n    @(\"356\u00d78=\", \"614\u00d76=\"),\n    @(\"995\u00d72=\", \"138\u00d77=\"),\n    @(\"534\u00d75=\", \"452\u00d77=\"),\n    @(\"329\u00d77=\", \"806\u00d72=\"),\n    @(\"196\u00d75=\", \"188\u00d78=\"),\n    @(\"949\u00d79=\", \"394\u00d77=\"),\n    @(\"921\u00d74=\", \"765\u00d73=\"),\n    @(\"459\u00d72=\", \"832\u00d77=\"),\n    @(\"242\u00d79=\", \"502\u00d75=\"),\n    @(\"749\u00d72=\", \"553\u00d72=\"),\n    @(\"692\u00d77=\", \"782\u00d75=\"),\n    @(\"707\u00d78=\", \"718\u00d73=\"),\n    @(\"315\u00d78=\", \"923\u00d78=\"),\n    @(\"515\u00d76=\", \"892\u00d77=\"),\n    @(\"398\u00d73=\", \"548\u00d77=\"),\n    @(\"325\u00d77=\", \"975\u00d78=\"),\n    @(\"465\u00d72=\", \"929\u00d78=\"),\n    @(\"353\u00d78=\", \"641\u00d75=\"),\n    @(\"647\u00d72=\", \"325\u00d76=\"),\n    @(\"963\u00d72=\", \"336\u00d74=\"),\n    @(\"422\u00d74=\", \"776\u00d77=\"),\n    @(\"645\u00d79=\", \"586\u00d76=\"),\n    @(\"362\u00d76=\", \"489\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
